$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block (rows 364-428) is shifted down by 2 rows: two new weekly
# observations are inserted at the top of the block (new rows 364-365), the
# remaining historical rows shift down by two positions, and the two rows
# that fall off the bottom of the original range become new rows 429-430.
# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T are constant for every row in this block,
# so only D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg) need to
# be rewritten for rows 364-428; rows 429-430 are brand new rows and need the
# constant columns written as well.

# Row 364
$ws.Range("D364").Value = 44522
$ws.Range("D364").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L364").Value = 'Pintón'
$ws.Range("M364").Value = 500
$ws.Range("N364").Value = 19000
$ws.Range("O364").Value = 19000
$ws.Range("P364").Value = 19000
$ws.Range("S364").Value = 950

# Row 365
$ws.Range("D365").Value = 44522
$ws.Range("D365").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L365").Value = 'Primera Pintón'
$ws.Range("M365").Value = 380
$ws.Range("N365").Value = 20000
$ws.Range("O365").Value = 20000
$ws.Range("P365").Value = 20000
$ws.Range("S365").Value = 1000

# Row 366
$ws.Range("D366").Value = 44246
$ws.Range("D366").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L366").Value = 'Pintón'
$ws.Range("M366").Value = 300
$ws.Range("N366").Value = 9000
$ws.Range("O366").Value = 9000
$ws.Range("P366").Value = 9000
$ws.Range("S366").Value = 450

# Row 367
$ws.Range("D367").Value = 44246
$ws.Range("D367").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L367").Value = 'Primera Pintón'
$ws.Range("M367").Value = 400
$ws.Range("N367").Value = 10000
$ws.Range("O367").Value = 10000
$ws.Range("P367").Value = 10000
$ws.Range("S367").Value = 500

# Row 368
$ws.Range("D368").Value = 44491
$ws.Range("D368").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L368").Value = 'Pintón'
$ws.Range("M368").Value = 1140
$ws.Range("N368").Value = 22000
$ws.Range("O368").Value = 23000
$ws.Range("P368").Value = 22526
$ws.Range("S368").Value = 1126

# Row 369
$ws.Range("D369").Value = 44272
$ws.Range("D369").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L369").Value = 'Pintón'
$ws.Range("M369").Value = 800
$ws.Range("N369").Value = 11000
$ws.Range("O369").Value = 11000
$ws.Range("P369").Value = 11000
$ws.Range("S369").Value = 550

# Row 370
$ws.Range("D370").Value = 44272
$ws.Range("D370").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L370").Value = 'Primera Pintón'
$ws.Range("M370").Value = 260
$ws.Range("N370").Value = 12000
$ws.Range("O370").Value = 12000
$ws.Range("P370").Value = 12000
$ws.Range("S370").Value = 600

# Row 371
$ws.Range("D371").Value = 44305
$ws.Range("D371").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L371").Value = 'Pintón'
$ws.Range("M371").Value = 1500
$ws.Range("N371").Value = 10000
$ws.Range("O371").Value = 10000
$ws.Range("P371").Value = 10000
$ws.Range("S371").Value = 500

# Row 372
$ws.Range("D372").Value = 44305
$ws.Range("D372").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L372").Value = 'Primera Pintón'
$ws.Range("M372").Value = 500
$ws.Range("N372").Value = 12000
$ws.Range("O372").Value = 12000
$ws.Range("P372").Value = 12000
$ws.Range("S372").Value = 600

# Row 373
$ws.Range("D373").Value = 44166
$ws.Range("D373").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L373").Value = 'Pintón'
$ws.Range("M373").Value = 400
$ws.Range("N373").Value = 16000
$ws.Range("O373").Value = 16000
$ws.Range("P373").Value = 16000
$ws.Range("S373").Value = 800

# Row 374
$ws.Range("D374").Value = 44166
$ws.Range("D374").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L374").Value = 'Primera Pintón'
$ws.Range("M374").Value = 500
$ws.Range("N374").Value = 17000
$ws.Range("O374").Value = 17000
$ws.Range("P374").Value = 17000
$ws.Range("S374").Value = 850

# Row 375
$ws.Range("D375").Value = 44225
$ws.Range("D375").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L375").Value = 'Pintón'
$ws.Range("M375").Value = 850
$ws.Range("N375").Value = 13000
$ws.Range("O375").Value = 13000
$ws.Range("P375").Value = 13000
$ws.Range("S375").Value = 650

# Row 376
$ws.Range("D376").Value = 44225
$ws.Range("D376").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L376").Value = 'Primera Pintón'
$ws.Range("M376").Value = 500
$ws.Range("N376").Value = 14000
$ws.Range("O376").Value = 14000
$ws.Range("P376").Value = 14000
$ws.Range("S376").Value = 700

# Row 377
$ws.Range("D377").Value = 44447
$ws.Range("D377").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L377").Value = 'Pintón'
$ws.Range("M377").Value = 300
$ws.Range("N377").Value = 19000
$ws.Range("O377").Value = 19000
$ws.Range("P377").Value = 19000
$ws.Range("S377").Value = 950

# Row 378
$ws.Range("D378").Value = 44447
$ws.Range("D378").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L378").Value = 'Primera Pintón'
$ws.Range("M378").Value = 400
$ws.Range("N378").Value = 20000
$ws.Range("O378").Value = 20000
$ws.Range("P378").Value = 20000
$ws.Range("S378").Value = 1000

# Row 379
$ws.Range("D379").Value = 44425
$ws.Range("D379").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L379").Value = 'Pintón'
$ws.Range("M379").Value = 1050
$ws.Range("N379").Value = 14000
$ws.Range("O379").Value = 14000
$ws.Range("P379").Value = 14000
$ws.Range("S379").Value = 700

# Row 380
$ws.Range("D380").Value = 44425
$ws.Range("D380").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L380").Value = 'Primera Pintón'
$ws.Range("M380").Value = 800
$ws.Range("N380").Value = 15000
$ws.Range("O380").Value = 15000
$ws.Range("P380").Value = 15000
$ws.Range("S380").Value = 750

# Row 381
$ws.Range("D381").Value = 44315
$ws.Range("D381").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L381").Value = 'Pintón'
$ws.Range("M381").Value = 850
$ws.Range("N381").Value = 13000
$ws.Range("O381").Value = 13000
$ws.Range("P381").Value = 13000
$ws.Range("S381").Value = 650

# Row 382
$ws.Range("D382").Value = 44315
$ws.Range("D382").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L382").Value = 'Primera Pintón'
$ws.Range("M382").Value = 400
$ws.Range("N382").Value = 15000
$ws.Range("O382").Value = 15000
$ws.Range("P382").Value = 15000
$ws.Range("S382").Value = 750

# Row 383
$ws.Range("D383").Value = 44348
$ws.Range("D383").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L383").Value = 'Maduro'
$ws.Range("M383").Value = 260
$ws.Range("N383").Value = 8000
$ws.Range("O383").Value = 8000
$ws.Range("P383").Value = 8000
$ws.Range("S383").Value = 400

# Row 384
$ws.Range("D384").Value = 44348
$ws.Range("D384").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L384").Value = 'Pintón'
$ws.Range("M384").Value = 850
$ws.Range("N384").Value = 10000
$ws.Range("O384").Value = 10000
$ws.Range("P384").Value = 10000
$ws.Range("S384").Value = 500

# Row 385
$ws.Range("D385").Value = 44348
$ws.Range("D385").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L385").Value = 'Primera Pintón'
$ws.Range("M385").Value = 600
$ws.Range("N385").Value = 12000
$ws.Range("O385").Value = 12000
$ws.Range("P385").Value = 12000
$ws.Range("S385").Value = 600

# Row 386
$ws.Range("D386").Value = 44322
$ws.Range("D386").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L386").Value = 'Pintón'
$ws.Range("M386").Value = 700
$ws.Range("N386").Value = 15000
$ws.Range("O386").Value = 15000
$ws.Range("P386").Value = 15000
$ws.Range("S386").Value = 750

# Row 387
$ws.Range("D387").Value = 44322
$ws.Range("D387").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L387").Value = 'Primera Pintón'
$ws.Range("M387").Value = 300
$ws.Range("N387").Value = 16000
$ws.Range("O387").Value = 16000
$ws.Range("P387").Value = 16000
$ws.Range("S387").Value = 800

# Row 388
$ws.Range("D388").Value = 44495
$ws.Range("D388").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L388").Value = 'Primera Pintón'
$ws.Range("M388").Value = 600
$ws.Range("N388").Value = 21000
$ws.Range("O388").Value = 21000
$ws.Range("P388").Value = 21000
$ws.Range("S388").Value = 1050

# Row 389
$ws.Range("D389").Value = 44232
$ws.Range("D389").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L389").Value = 'Pintón'
$ws.Range("M389").Value = 840
$ws.Range("N389").Value = 13000
$ws.Range("O389").Value = 13000
$ws.Range("P389").Value = 13000
$ws.Range("S389").Value = 650

# Row 390
$ws.Range("D390").Value = 44232
$ws.Range("D390").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L390").Value = 'Primera Pintón'
$ws.Range("M390").Value = 310
$ws.Range("N390").Value = 16000
$ws.Range("O390").Value = 16000
$ws.Range("P390").Value = 16000
$ws.Range("S390").Value = 800

# Row 391
$ws.Range("D391").Value = 44327
$ws.Range("D391").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L391").Value = 'Pintón'
$ws.Range("M391").Value = 1000
$ws.Range("N391").Value = 14000
$ws.Range("O391").Value = 14000
$ws.Range("P391").Value = 14000
$ws.Range("S391").Value = 700

# Row 392
$ws.Range("D392").Value = 44510
$ws.Range("D392").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L392").Value = 'Pintón'
$ws.Range("M392").Value = 800
$ws.Range("N392").Value = 15000
$ws.Range("O392").Value = 15000
$ws.Range("P392").Value = 15000
$ws.Range("S392").Value = 750

# Row 393
$ws.Range("D393").Value = 44510
$ws.Range("D393").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L393").Value = 'Primera Pintón'
$ws.Range("M393").Value = 500
$ws.Range("N393").Value = 17000
$ws.Range("O393").Value = 17000
$ws.Range("P393").Value = 17000
$ws.Range("S393").Value = 850

# Row 394
$ws.Range("D394").Value = 44161
$ws.Range("D394").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L394").Value = 'Pintón'
$ws.Range("M394").Value = 1400
$ws.Range("N394").Value = 14000
$ws.Range("O394").Value = 14000
$ws.Range("P394").Value = 14000
$ws.Range("S394").Value = 700

# Row 395
$ws.Range("D395").Value = 44161
$ws.Range("D395").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L395").Value = 'Primera Pintón'
$ws.Range("M395").Value = 500
$ws.Range("N395").Value = 15000
$ws.Range("O395").Value = 15000
$ws.Range("P395").Value = 15000
$ws.Range("S395").Value = 750

# Row 396
$ws.Range("D396").Value = 44468
$ws.Range("D396").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L396").Value = 'Pintón'
$ws.Range("M396").Value = 800
$ws.Range("N396").Value = 12000
$ws.Range("O396").Value = 12000
$ws.Range("P396").Value = 12000
$ws.Range("S396").Value = 600

# Row 397
$ws.Range("D397").Value = 44468
$ws.Range("D397").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L397").Value = 'Primera Pintón'
$ws.Range("M397").Value = 450
$ws.Range("N397").Value = 14000
$ws.Range("O397").Value = 14000
$ws.Range("P397").Value = 14000
$ws.Range("S397").Value = 700

# Row 398
$ws.Range("D398").Value = 44517
$ws.Range("D398").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L398").Value = 'Pintón'
$ws.Range("M398").Value = 12000
$ws.Range("N398").Value = 14000
$ws.Range("O398").Value = 14000
$ws.Range("P398").Value = 14000
$ws.Range("S398").Value = 700

# Row 399
$ws.Range("D399").Value = 44517
$ws.Range("D399").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L399").Value = 'Primera Pintón'
$ws.Range("M399").Value = 600
$ws.Range("N399").Value = 15000
$ws.Range("O399").Value = 15000
$ws.Range("P399").Value = 15000
$ws.Range("S399").Value = 750

# Row 400
$ws.Range("D400").Value = 44238
$ws.Range("D400").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L400").Value = 'Pintón'
$ws.Range("M400").Value = 400
$ws.Range("N400").Value = 7000
$ws.Range("O400").Value = 7000
$ws.Range("P400").Value = 7000
$ws.Range("S400").Value = 350

# Row 401
$ws.Range("D401").Value = 44238
$ws.Range("D401").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L401").Value = 'Primera Pintón'
$ws.Range("M401").Value = 600
$ws.Range("N401").Value = 8000
$ws.Range("O401").Value = 8000
$ws.Range("P401").Value = 8000
$ws.Range("S401").Value = 400

# Row 402
$ws.Range("D402").Value = 44391
$ws.Range("D402").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L402").Value = 'Pintón'
$ws.Range("M402").Value = 300
$ws.Range("N402").Value = 9000
$ws.Range("O402").Value = 9000
$ws.Range("P402").Value = 9000
$ws.Range("S402").Value = 450

# Row 403
$ws.Range("D403").Value = 44391
$ws.Range("D403").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L403").Value = 'Primera Pintón'
$ws.Range("M403").Value = 500
$ws.Range("N403").Value = 11000
$ws.Range("O403").Value = 11000
$ws.Range("P403").Value = 11000
$ws.Range("S403").Value = 550

# Row 404
$ws.Range("D404").Value = 44389
$ws.Range("D404").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L404").Value = 'Pintón'
$ws.Range("M404").Value = 300
$ws.Range("N404").Value = 10000
$ws.Range("O404").Value = 10000
$ws.Range("P404").Value = 10000
$ws.Range("S404").Value = 500

# Row 405
$ws.Range("D405").Value = 44389
$ws.Range("D405").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L405").Value = 'Primera Pintón'
$ws.Range("M405").Value = 450
$ws.Range("N405").Value = 12000
$ws.Range("O405").Value = 12000
$ws.Range("P405").Value = 12000
$ws.Range("S405").Value = 600

# Row 406
$ws.Range("D406").Value = 44251
$ws.Range("D406").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L406").Value = 'Pintón'
$ws.Range("M406").Value = 560
$ws.Range("N406").Value = 9000
$ws.Range("O406").Value = 9000
$ws.Range("P406").Value = 9000
$ws.Range("S406").Value = 450

# Row 407
$ws.Range("D407").Value = 44251
$ws.Range("D407").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L407").Value = 'Primera Pintón'
$ws.Range("M407").Value = 340
$ws.Range("N407").Value = 10000
$ws.Range("O407").Value = 10000
$ws.Range("P407").Value = 10000
$ws.Range("S407").Value = 500

# Row 408
$ws.Range("D408").Value = 44340
$ws.Range("D408").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L408").Value = 'Pintón'
$ws.Range("M408").Value = 1250
$ws.Range("N408").Value = 10000
$ws.Range("O408").Value = 10000
$ws.Range("P408").Value = 10000
$ws.Range("S408").Value = 500

# Row 409
$ws.Range("D409").Value = 44340
$ws.Range("D409").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L409").Value = 'Primera Pintón'
$ws.Range("M409").Value = 480
$ws.Range("N409").Value = 11500
$ws.Range("O409").Value = 11500
$ws.Range("P409").Value = 11500
$ws.Range("S409").Value = 575

# Row 410
$ws.Range("D410").Value = 44515
$ws.Range("D410").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L410").Value = 'Pintón'
$ws.Range("M410").Value = 1000
$ws.Range("N410").Value = 15000
$ws.Range("O410").Value = 15000
$ws.Range("P410").Value = 15000
$ws.Range("S410").Value = 750

# Row 411
$ws.Range("D411").Value = 44515
$ws.Range("D411").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L411").Value = 'Primera Pintón'
$ws.Range("M411").Value = 500
$ws.Range("N411").Value = 17000
$ws.Range("O411").Value = 17000
$ws.Range("P411").Value = 17000
$ws.Range("S411").Value = 850

# Row 412
$ws.Range("D412").Value = 44330
$ws.Range("D412").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L412").Value = 'Pintón'
$ws.Range("M412").Value = 800
$ws.Range("N412").Value = 10000
$ws.Range("O412").Value = 10000
$ws.Range("P412").Value = 10000
$ws.Range("S412").Value = 500

# Row 413
$ws.Range("D413").Value = 44330
$ws.Range("D413").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L413").Value = 'Primera Pintón'
$ws.Range("M413").Value = 500
$ws.Range("N413").Value = 12000
$ws.Range("O413").Value = 12000
$ws.Range("P413").Value = 12000
$ws.Range("S413").Value = 600

# Row 414
$ws.Range("D414").Value = 44432
$ws.Range("D414").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L414").Value = 'Pintón'
$ws.Range("M414").Value = 500
$ws.Range("N414").Value = 12000
$ws.Range("O414").Value = 12000
$ws.Range("P414").Value = 12000
$ws.Range("S414").Value = 600

# Row 415
$ws.Range("D415").Value = 44432
$ws.Range("D415").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L415").Value = 'Primera Pintón'
$ws.Range("M415").Value = 500
$ws.Range("N415").Value = 13000
$ws.Range("O415").Value = 13000
$ws.Range("P415").Value = 13000
$ws.Range("S415").Value = 650

# Row 416
$ws.Range("D416").Value = 44181
$ws.Range("D416").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L416").Value = 'Pintón'
$ws.Range("M416").Value = 600
$ws.Range("N416").Value = 12000
$ws.Range("O416").Value = 12000
$ws.Range("P416").Value = 12000
$ws.Range("S416").Value = 600

# Row 417
$ws.Range("D417").Value = 44181
$ws.Range("D417").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L417").Value = 'Primera Pintón'
$ws.Range("M417").Value = 300
$ws.Range("N417").Value = 13000
$ws.Range("O417").Value = 13000
$ws.Range("P417").Value = 13000
$ws.Range("S417").Value = 650

# Row 418
$ws.Range("D418").Value = 44194
$ws.Range("D418").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L418").Value = 'Pintón'
$ws.Range("M418").Value = 300
$ws.Range("N418").Value = 13000
$ws.Range("O418").Value = 13000
$ws.Range("P418").Value = 13000
$ws.Range("S418").Value = 650

# Row 419
$ws.Range("D419").Value = 44194
$ws.Range("D419").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L419").Value = 'Primera Pintón'
$ws.Range("M419").Value = 500
$ws.Range("N419").Value = 12000
$ws.Range("O419").Value = 12000
$ws.Range("P419").Value = 12000
$ws.Range("S419").Value = 600

# Row 420
$ws.Range("D420").Value = 44271
$ws.Range("D420").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L420").Value = 'Pintón'
$ws.Range("M420").Value = 800
$ws.Range("N420").Value = 12000
$ws.Range("O420").Value = 12000
$ws.Range("P420").Value = 12000
$ws.Range("S420").Value = 600

# Row 421
$ws.Range("D421").Value = 44271
$ws.Range("D421").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L421").Value = 'Primera Pintón'
$ws.Range("M421").Value = 360
$ws.Range("N421").Value = 13000
$ws.Range("O421").Value = 13000
$ws.Range("P421").Value = 13000
$ws.Range("S421").Value = 650

# Row 422
$ws.Range("D422").Value = 44307
$ws.Range("D422").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L422").Value = 'Pintón'
$ws.Range("M422").Value = 600
$ws.Range("N422").Value = 11000
$ws.Range("O422").Value = 11000
$ws.Range("P422").Value = 11000
$ws.Range("S422").Value = 550

# Row 423
$ws.Range("D423").Value = 44307
$ws.Range("D423").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L423").Value = 'Primera Pintón'
$ws.Range("M423").Value = 300
$ws.Range("N423").Value = 13000
$ws.Range("O423").Value = 13000
$ws.Range("P423").Value = 13000
$ws.Range("S423").Value = 650

# Row 424
$ws.Range("D424").Value = 44400
$ws.Range("D424").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L424").Value = 'Pintón'
$ws.Range("M424").Value = 150
$ws.Range("N424").Value = 18000
$ws.Range("O424").Value = 18000
$ws.Range("P424").Value = 18000
$ws.Range("S424").Value = 900

# Row 425
$ws.Range("D425").Value = 44400
$ws.Range("D425").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L425").Value = 'Primera Pintón'
$ws.Range("M425").Value = 100
$ws.Range("N425").Value = 20000
$ws.Range("O425").Value = 20000
$ws.Range("P425").Value = 20000
$ws.Range("S425").Value = 1000

# Row 426
$ws.Range("D426").Value = 44309
$ws.Range("D426").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L426").Value = 'Pintón'
$ws.Range("M426").Value = 800
$ws.Range("N426").Value = 10000
$ws.Range("O426").Value = 10000
$ws.Range("P426").Value = 10000
$ws.Range("S426").Value = 500

# Row 427
$ws.Range("D427").Value = 44309
$ws.Range("D427").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L427").Value = 'Primera Pintón'
$ws.Range("M427").Value = 600
$ws.Range("N427").Value = 12000
$ws.Range("O427").Value = 12000
$ws.Range("P427").Value = 12000
$ws.Range("S427").Value = 600

# Row 428
$ws.Range("D428").Value = 44508
$ws.Range("D428").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L428").Value = 'Pintón'
$ws.Range("M428").Value = 1050
$ws.Range("N428").Value = 15000
$ws.Range("O428").Value = 15000
$ws.Range("P428").Value = 15000
$ws.Range("S428").Value = 750

# Row 429
$ws.Range("A429").Value = 5
$ws.Range("B429").Value = 'Macroferia Regional de Talca'
$ws.Range("C429").Value = 'Maule'
$ws.Range("E429").Value = 7
$ws.Range("F429").Value = 'Fruta'
$ws.Range("G429").Value = 100108
$ws.Range("H429").Value = 'Tropicales y subtropicales'
$ws.Range("I429").Value = 100108006
$ws.Range("J429").Value = 'Plátano'
$ws.Range("K429").Value = 'Sin especificar'
$ws.Range("Q429").Value = '$/caja 20 kilos'
$ws.Range("R429").Value = 'Ecuador'
$ws.Range("T429").Value = 20
$ws.Range("D429").Value = 44508
$ws.Range("D429").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L429").Value = 'Primera Pintón'
$ws.Range("M429").Value = 540
$ws.Range("N429").Value = 17000
$ws.Range("O429").Value = 17000
$ws.Range("P429").Value = 17000
$ws.Range("S429").Value = 850

# Row 430
$ws.Range("A430").Value = 5
$ws.Range("B430").Value = 'Macroferia Regional de Talca'
$ws.Range("C430").Value = 'Maule'
$ws.Range("E430").Value = 7
$ws.Range("F430").Value = 'Fruta'
$ws.Range("G430").Value = 100108
$ws.Range("H430").Value = 'Tropicales y subtropicales'
$ws.Range("I430").Value = 100108006
$ws.Range("J430").Value = 'Plátano'
$ws.Range("K430").Value = 'Sin especificar'
$ws.Range("Q430").Value = '$/caja 20 kilos'
$ws.Range("R430").Value = 'Ecuador'
$ws.Range("T430").Value = 20
$ws.Range("D430").Value = 44201
$ws.Range("D430").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("L430").Value = 'Pintón'
$ws.Range("M430").Value = 800
$ws.Range("N430").Value = 15000
$ws.Range("O430").Value = 15000
$ws.Range("P430").Value = 15000
$ws.Range("S430").Value = 750

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()